$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '60.866.27'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -2.47%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.394.69'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -2.41%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.09%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '570.62'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.50%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '139.74'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.69%  '
$ws.Range("E7").Value = '  -0.04%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.527'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.37%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.393.45'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -2.28%  '
$ws.Range("E10").Value = '  -1.00%  '
$ws.Range("E11").Value = '  -0.27%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '5.09'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -2.31%  '
$ws.Range("E13").Value = '  -2.33%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '25.95'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.86%  '
$ws.Range("E15").Value = '  -3.85%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.844.64'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.84%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '60.664.38'
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.403.16'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.91%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.48'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +5.01%  '
$ws.Range("E20").Value = '  -2.02%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '322.49'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.89%  '
$ws.Range("E22").Value = '  -1.79%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.04'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.71%  '
$ws.Range("E24").Value = '  +0.06%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.87'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -5.06%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '64.78'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -1.32%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '583.26'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.34%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.42'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -9.62%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.516.20'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -2.33%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0₃0921'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -4.79%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.92'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.57%  '
$ws.Range("E32").Value = '  -5.95%  '
$ws.Range("E33").Value = '  -2.23%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.133'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -2.92%  '
$ws.Range("E35").Value = '  -0.26%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.64'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -6.09%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.41'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -3.15%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '151.18'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -2.10%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.368'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -2.99%  '
$ws.Range("E40").Value = '  -0.87%  '
$ws.Range("E41").Value = '  -3.85%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.999'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.01%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.68'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -3.15%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '41.16'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -4.38%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.35'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -4.71%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0₆0287'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +12.35%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '140.81'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.60%  '
$ws.Range("E48").Value = '  -3.99%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.590'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -3.06%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '19.54'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.62%  '
$ws.Range("E51").Value = '  -3.72%  '
